$wb = $excel.ActiveWorkbook

# Each entry: sheet name, row number, then a map of column letter -> new value ($null means clear the cell)
$edits = @(
    @{ Sheet = "ALC"; Row = 116; Cols = @{ "H" = 3685.625; "I" = 4201; "J" = 2826.6667; "K" = 4201; "L" = 2826.6667; "M" = -759; "N" = -9710.6667 } }
    @{ Sheet = "ALC"; Row = 117; Cols = @{ "H" = 48000; "J" = 48000; "L" = 48000; "N" = -57178 } }
    @{ Sheet = "ALC"; Row = 118; Cols = @{ "H" = 1017.0732; "I" = 482.30768; "J" = 1265.3572; "K" = 1446.92304; "L" = 3796.0716; "M" = 210.0769599999999; "N" = -7110.071599999999 } }
    @{ Sheet = "ALC"; Row = 120; Cols = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null } }
    @{ Sheet = "ARM"; Row = 2; Cols = @{ "H" = 111965.555; "I" = 250672.75; "J" = 999.8; "K" = 250672.75; "L" = 999.8; "M" = -250559.75; "N" = -1225.8 } }
    @{ Sheet = "ARM"; Row = 36; Cols = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "ARM"; Row = 88; Cols = @{ "H" = 3792.4443; "I" = 3140.2856; "J" = 6075; "K" = 3140.2856; "L" = 6075; "M" = -2734.2856; "N" = -6887 } }
    @{ Sheet = "ARM"; Row = 91; Cols = @{ "H" = 3792.4443; "I" = 3140.2856; "J" = 6075; "K" = 3140.2856; "L" = 6075; "M" = -1736.2856; "N" = -8883 } }
    @{ Sheet = "ARM"; Row = 102; Cols = @{ "H" = 1014.6; "I" = 944.2143; "J" = 2000; "K" = 944.2143; "L" = 2000; "M" = 677.7857; "N" = -5244 } }
    @{ Sheet = "ARM"; Row = 116; Cols = @{ "H" = 111965.555; "I" = 250672.75; "J" = 999.8; "K" = 250672.75; "L" = 999.8; "M" = -248378.75; "N" = -5587.8 } }
    @{ Sheet = "ARM"; Row = 117; Cols = @{ "H" = 32478.6; "J" = 32478.6; "L" = 32478.6; "N" = -41656.6 } }
    @{ Sheet = "ARM"; Row = 132; Cols = @{ "H" = 4980.868; "I" = 3773.3057; "K" = 11319.9171; "M" = -8789.917099999999 } }
    @{ Sheet = "BSM"; Row = 3; Cols = @{ "H" = 111965.555; "I" = 250672.75; "J" = 999.8; "K" = 250672.75; "L" = 999.8; "M" = -250558.75; "N" = -1227.8 } }
    @{ Sheet = "BSM"; Row = 33; Cols = @{ "H" = 2020.6666; "I" = 2020.6666; "K" = 2020.6666; "M" = -1684.6666 } }
    @{ Sheet = "BSM"; Row = 86; Cols = @{ "H" = 2096.3794; "I" = 1848.1333; "J" = 2362.3572; "K" = 1848.1333; "L" = 2362.3572; "M" = -725.1333; "N" = -4608.3572 } }
    @{ Sheet = "BSM"; Row = 89; Cols = @{ "H" = 2096.3794; "I" = 1848.1333; "J" = 2362.3572; "K" = 9240.666499999999; "L" = 11811.786; "M" = -3624.666499999999; "N" = -23043.786 } }
    @{ Sheet = "BSM"; Row = 116; Cols = @{ "H" = 48000; "J" = 48000; "L" = 48000; "N" = -57178 } }
    @{ Sheet = "BSM"; Row = 117; Cols = @{ "H" = 48000; "J" = 48000; "L" = 48000; "N" = -57178 } }
    @{ Sheet = "BSM"; Row = 118; Cols = @{ "H" = 48000; "J" = 48000; "L" = 48000; "N" = -51314 } }
    @{ Sheet = "BSM"; Row = 134; Cols = @{ "H" = 729.725; "I" = 701.52; "J" = 1152.8; "K" = 2104.56; "L" = 3458.4; "M" = 430.4400000000001; "N" = -8528.4 } }
    @{ Sheet = "CRP"; Row = 39; Cols = @{ "H" = 5000; "I" = 3000; "K" = 3000; "M" = -2609 } }
    @{ Sheet = "CRP"; Row = 49; Cols = @{ "H" = 5000; "I" = 3000; "K" = 3000; "M" = -2818 } }
    @{ Sheet = "CRP"; Row = 58; Cols = @{ "H" = 841.9216; "J" = 1275.4286; "L" = 1275.4286; "N" = -1681.4286 } }
    @{ Sheet = "CRP"; Row = 105; Cols = @{ "H" = 1395; "I" = 1012.1429; "J" = 1730; "K" = 1012.1429; "L" = 1730; "M" = 734.8570999999999; "N" = -5224 } }
    @{ Sheet = "CRP"; Row = 132; Cols = @{ "H" = 19612018; "I" = 29416764; "J" = 2524.353; "K" = 88250292; "L" = 7573.059; "M" = -88247762; "N" = -12633.059 } }
    @{ Sheet = "CRP"; Row = 136; Cols = @{ "H" = 841.9216; "J" = 1275.4286; "L" = 3826.2858; "N" = -8926.2858 } }
    @{ Sheet = "CUL"; Row = 117; Cols = @{ "H" = 1097.4286; "I" = 966.6667; "J" = 1195.5; "K" = 2900.0001; "L" = 3586.5; "M" = 541.9998999999998; "N" = -10470.5 } }
    @{ Sheet = "CUL"; Row = 119; Cols = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null } }
    @{ Sheet = "CUL"; Row = 120; Cols = @{ "H" = 7166.5; "I" = 1000; "J" = 13333; "K" = 3000; "L" = 39999; "M" = 1838; "N" = -49675 } }
    @{ Sheet = "CUL"; Row = 131; Cols = @{ "H" = 600.90625; "J" = 996; "L" = 2988; "N" = -13068 } }
    @{ Sheet = "GSM"; Row = 116; Cols = @{ "H" = 35001; "J" = 35001; "L" = 35001; "N" = -44179 } }
    @{ Sheet = "GSM"; Row = 117; Cols = @{ "H" = 9200.333000000001; "J" = 9200.333000000001; "L" = 9200.333000000001; "N" = -16084.333 } }
    @{ Sheet = "GSM"; Row = 118; Cols = @{ "H" = 17555; "I" = 0; "J" = 17555; "K" = 0; "L" = 17555; "M" = $null; "N" = -20869 } }
    @{ Sheet = "LTW"; Row = 132; Cols = @{ "H" = 4264.019; "I" = 5824.5356; "J" = 2443.4167; "K" = 17473.6068; "L" = 7330.250100000001; "M" = -14943.6068; "N" = -12390.2501 } }
    @{ Sheet = "WVR"; Row = 5; Cols = @{ "H" = 774769.25; "I" = 3000; "J" = 1006300; "K" = 3000; "L" = 1006300; "M" = -2888; "N" = -1006524 } }
    @{ Sheet = "WVR"; Row = 117; Cols = @{ "H" = 47999; "J" = 47999; "L" = 47999; "N" = -57177 } }
    @{ Sheet = "WVR"; Row = 121; Cols = @{ "H" = 48000; "J" = 48000; "L" = 48000; "N" = -51494 } }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($col in $edit.Cols.Keys) {
        $addr = "$col$($edit.Row)"
        $val = $edit.Cols[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value2 = $val
        }
    }
}

Write-Host "Applied $($edits.Count) row edits."
